$p = $ppt.ActivePresentation

# Slide 5 (the closing / "thank you" slide) gets a new textbox with a
# link to the rpubs.com write-up of the analysis.
$s = $p.Slides.Item(5)

# Coordinates/size come straight from the target OOXML (EMU); the
# PowerPoint object model's Shapes.AddTextbox takes points, so convert
# (1 pt = 12700 EMU).
$left   = 2391508 / 12700
$top    = 2433711 / 12700
$width  = 4139275 / 12700
$height = 369332 / 12700

$tb = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
$tb.Name = "TextBox 2"

# Matches a plain "click and type" textbox: no fill, shrink-to-fit text
# (no wrapping, box auto-sizes to the typed line).
$tb.Fill.Visible = $false
$tb.TextFrame.WordWrap = $false
$tb.TextFrame.AutoSize = 1

$tb.TextFrame.TextRange.Text = "http://rpubs.com/azureblue83/327856"

Write-Output "done"
